$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- Update time_taken (column F) timestamps for rows 2..96 on the data sheet ---
$newTimestamps = @(
    "2021-10-05 14:33:56.422672",
    "2021-10-05 14:33:56.422680",
    "2021-10-05 14:33:56.422683",
    "2021-10-05 14:33:56.422686",
    "2021-10-05 14:33:56.422688",
    "2021-10-05 14:33:56.422691",
    "2021-10-05 14:33:56.422694",
    "2021-10-05 14:33:56.422696",
    "2021-10-05 14:33:56.422699",
    "2021-10-05 14:33:56.422701",
    "2021-10-05 14:33:56.422704",
    "2021-10-05 14:33:56.422706",
    "2021-10-05 14:33:56.422709",
    "2021-10-05 14:33:56.422711",
    "2021-10-05 14:33:56.422714",
    "2021-10-05 14:33:56.422716",
    "2021-10-05 14:33:56.422719",
    "2021-10-05 14:33:56.422722",
    "2021-10-05 14:33:56.422724",
    "2021-10-05 14:33:56.422727",
    "2021-10-05 14:33:56.422729",
    "2021-10-05 14:33:56.422732",
    "2021-10-05 14:33:56.422734",
    "2021-10-05 14:33:56.422736",
    "2021-10-05 14:33:56.422739",
    "2021-10-05 14:33:56.422742",
    "2021-10-05 14:33:56.422744",
    "2021-10-05 14:33:56.422746",
    "2021-10-05 14:33:56.422749",
    "2021-10-05 14:33:56.422751",
    "2021-10-05 14:33:56.422754",
    "2021-10-05 14:33:56.422756",
    "2021-10-05 14:33:56.422759",
    "2021-10-05 14:33:56.422762",
    "2021-10-05 14:33:56.422764",
    "2021-10-05 14:33:56.422766",
    "2021-10-05 14:33:56.422769",
    "2021-10-05 14:33:56.422771",
    "2021-10-05 14:33:56.422774",
    "2021-10-05 14:33:56.422776",
    "2021-10-05 14:33:56.422779",
    "2021-10-05 14:33:56.422782",
    "2021-10-05 14:33:56.422784",
    "2021-10-05 14:33:56.422787",
    "2021-10-05 14:33:56.422789",
    "2021-10-05 14:33:56.422791",
    "2021-10-05 14:33:56.422794",
    "2021-10-05 14:33:56.422796",
    "2021-10-05 14:33:56.422799",
    "2021-10-05 14:33:56.422801",
    "2021-10-05 14:33:56.422804",
    "2021-10-05 14:33:56.422806",
    "2021-10-05 14:33:56.422809",
    "2021-10-05 14:33:56.422812",
    "2021-10-05 14:33:56.422814",
    "2021-10-05 14:33:56.422817",
    "2021-10-05 14:33:56.422819",
    "2021-10-05 14:33:56.422822",
    "2021-10-05 14:33:56.422824",
    "2021-10-05 14:33:56.422827",
    "2021-10-05 14:33:56.422829",
    "2021-10-05 14:33:56.422832",
    "2021-10-05 14:33:56.422834",
    "2021-10-05 14:33:56.422837",
    "2021-10-05 14:33:56.422840",
    "2021-10-05 14:33:56.422843",
    "2021-10-05 14:33:56.422845",
    "2021-10-05 14:33:56.422848",
    "2021-10-05 14:33:56.422850",
    "2021-10-05 14:33:56.422853",
    "2021-10-05 14:33:56.422855",
    "2021-10-05 14:33:56.422857",
    "2021-10-05 14:33:56.422860",
    "2021-10-05 14:33:56.422862",
    "2021-10-05 14:33:56.422865",
    "2021-10-05 14:33:56.422867",
    "2021-10-05 14:33:56.422872",
    "2021-10-05 14:33:56.422875",
    "2021-10-05 14:33:56.422877",
    "2021-10-05 14:33:56.422880",
    "2021-10-05 14:33:56.422883",
    "2021-10-05 14:33:56.422885",
    "2021-10-05 14:33:56.422887",
    "2021-10-05 14:33:56.422890",
    "2021-10-05 14:33:56.422892",
    "2021-10-05 14:33:56.422895",
    "2021-10-05 14:33:56.422897",
    "2021-10-05 14:33:56.422900",
    "2021-10-05 14:33:56.422902",
    "2021-10-05 14:33:56.422905",
    "2021-10-05 14:33:56.422907",
    "2021-10-05 14:33:56.422909",
    "2021-10-05 14:33:56.422913",
    "2021-10-05 14:33:56.422916",
    "2021-10-05 14:33:56.422918"
)

for ($i = 0; $i -lt $newTimestamps.Length; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 6).Value = $newTimestamps[$i]
}

# --- Add the new "metadata" worksheet, placed right after "data" ---
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

# Copy the header styling (bold font + border + centered) used on the data sheet
$dataSheet.Range("B1:F1").Copy()
$metaSheet.Range("B1:F1").PasteSpecial(-4122)
$dataSheet.Range("F1").Copy()
$metaSheet.Range("G1").PasteSpecial(-4122)

# Copy the index-column styling used for column A values
$dataSheet.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)

# Header row
$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

# Data row
$metaSheet.Range("A2").Value = 0
$metaSheet.Range("B2").Value = "Hand and foot malformations"
$metaSheet.Range("C2").Value = 3729
$metaSheet.Range("D2").NumberFormat = "@"
$metaSheet.Range("D2").Value = "0.55"
$metaSheet.Range("E2").Value = "2021-09-23T01:51:42.163707Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:33:56.419046"
$metaSheet.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/3729/?format=json"

Write-Host "edit complete"
